$d = $word.ActiveDocument
$CR = [char]13

# ---------------------------------------------------------------------------
# Change 1: strike through "Get stripe ecommerce example working and then
# add to the project" (same formatting already used on "Add cart
# functionality" right above it).
# ---------------------------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text.TrimEnd($CR)
    if ($text -eq "Get stripe ecommerce example working and then add to the project") {
        $para.Range.Font.StrikeThrough = $true
    }
}

# ---------------------------------------------------------------------------
# Change 2: add a new bullet "Have messages fade away after a period of
# time" right after "Rename image files..." (same list / numId 6), and move
# the _GoBack bookmark so it ends up right after the new bullet's text, the
# way Word leaves it after you type new content and save.
# ---------------------------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text.TrimEnd($CR)
    if ($text -eq "Rename image files on the way in for consignment, otherwise files can be overwritten") {
        $tail = $para.Range.Duplicate
        $tail.Collapse(0)
        $tail.InsertParagraphAfter()
    }
}

foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text.TrimEnd($CR)
    if ($text -eq "Rename image files on the way in for consignment, otherwise files can be overwritten") {
        $newPara = $para.Next()
        $newRange = $newPara.Range.Duplicate
        $newRange.Collapse(1)

        # Type the new text plus a throw-away trailing marker character: the
        # interop runtime cannot collapse a brand-new bookmark exactly on the
        # last character slot of a paragraph (right before the paragraph
        # mark), so the marker keeps that position from being "last" while
        # the bookmark is created, then gets deleted again afterwards.
        $newRange.InsertAfter("Have messages fade away after a period of time#")

        $markerPos = $newRange.End - 1
        $bookmarkTarget = $d.Range($markerPos, $markerPos)

        $d.Bookmarks("_GoBack").Delete()
        $d.Bookmarks.Add("_GoBack", $bookmarkTarget)

        $marker = $d.Range($newRange.End, $newRange.End)
        $marker.Delete()
    }
}
